$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new expense row 37: "Gastos de Parto" ---
$ws.Range("A37").Value = 150000
$ws.Range("B37").Value = "Gastos de Parto"
# Copy the date format from the row above (style only), then set the value
$ws.Range("C36").Copy()
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("C37").Value = 44854

# --- Add new expense row 38: "Gastos de alimento de engorde" ---
$ws.Range("A38").Value = 102000
$ws.Range("B38").Value = "Gastos de alimento de engorde"
$ws.Range("C36").Copy()
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("C38").Value = 44862

# --- Extend the table (Tabla1) to cover the two new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C38"))

# --- Update the total-gasto formula so it includes the new rows ---
$ws.Range("G16").Formula = "=SUM(A2:A38)"

# --- Update the selected/active cell like the saved workbook shows ---
$ws.Range("D38").Select()
